$d = $word.ActiveDocument

# Replace the first (and only expected) occurrence of $oldText in the document
# body with $newText. Uses Find (without its own Replace) + a direct Range.Text
# assignment so that straight quotes/apostrophes in the replacement text are not
# mangled by Word's "smart quotes" AutoFormat/AutoCorrect behaviour, and so the
# surrounding run's formatting (rPr/rStyle) is preserved.
function Replace-ExactText([string]$oldText, [string]$newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $newText
    } else {
        Write-Output "NOT FOUND: $oldText"
    }
    return $found
}

# 1. (Taxa,Region) -> (Taxa,Month,Region) (keep the preceding "group_by" run/style untouched)
Replace-ExactText "(Taxa,Region) " "(Taxa,Month,Region) "

# 2. summarise() grouped-output message: add 'Month' to the grouping columns mentioned
Replace-ExactText "## ``summarise()`` has grouped output by 'Taxa'. You can override using the" `
                   "## ``summarise()`` has grouped output by 'Taxa', 'Month'. You can override using the"

# 3. Joining with `by = join_by(Region, Taxa)` -> add Month
Replace-ExactText "## Joining with ``by = join_by(Region, Taxa)``" `
                   "## Joining with ``by = join_by(Month, Region, Taxa)``"

# 4. CSV filename now includes a date suffix
Replace-ExactText '"FlowZoopData_2022ROC_EffectsAnalysis_CohortYear.csv"' `
                   '"FlowZoopData_2022ROC_EffectsAnalysis_CohortYear_2024-09-05.csv"'
